# EMTBootCamp.xlsx - "updated the SCR ramping impedances"
#
# StepSCR sheet reworked from a single Xt/Zs/Rs impedance model to a
# Wind/Solar plant-reactance (Xplant) model, and the XS defined name
# (localSheetId=3, i.e. StepSCR) is repointed from the old H5 cell to
# the new J6 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StepSCR")

# Start from a clean sheet - the old layout (A1:H9) is being fully
# replaced by a new layout (A1:J12).
$ws.Cells.Clear()

# xlRight / xlLeft alignment constants (avoid relying on [Microsoft.Office...] enums)
$xlRight = -4152

# ---- Row 1: small two-column header above the Solar/Wind blocks ----
$c = $ws.Range("D1"); $c.Value = "Solar"; $c.Font.Bold = $true
$c = $ws.Range("F1"); $c.Value = "Wind";  $c.Font.Bold = $true

# ---- Row 2: main header row ----
$c = $ws.Range("A2"); $c.Value = "SCR";     $c.Font.Bold = $true; $c.HorizontalAlignment = $xlRight
$c = $ws.Range("B2"); $c.Value = "SCMVA";   $c.Font.Bold = $true; $c.HorizontalAlignment = $xlRight
$c = $ws.Range("C2"); $c.Value = "X1 [W]";  $c.Font.Bold = $true; $c.HorizontalAlignment = $xlRight
$c = $ws.Range("D2"); $c.Value = "Xplant";  $c.Font.Bold = $true; $c.HorizontalAlignment = $xlRight
$c = $ws.Range("E2"); $c.Value = "Xs [W]";  $c.Font.Bold = $true; $c.HorizontalAlignment = $xlRight
$c = $ws.Range("F2"); $c.Value = "Xplant";  $c.Font.Bold = $true; $c.HorizontalAlignment = $xlRight
$c = $ws.Range("G2"); $c.Value = "Xs [W]";  $c.Font.Bold = $true; $c.HorizontalAlignment = $xlRight
$c = $ws.Range("H2"); $c.Font.Bold = $true; $c.HorizontalAlignment = $xlRight
$c = $ws.Range("I2"); $c.Value = "kVs";     $c.Font.Bold = $true
$ws.Range("J2").Value = 230

$c = $ws.Range("I3"); $c.Value = "IBR MVA"; $c.Font.Bold = $true
$ws.Range("J3").Value = 100

# ---- Data rows 3..11 ----
$rows = @(
    @{r=3;  a=20},
    @{r=4;  a=10},
    @{r=5;  a=5},
    @{r=6;  a=4},
    @{r=7;  a=3},
    @{r=8;  a=2.5},
    @{r=9;  a=2},
    @{r=10; a=1.5},
    @{r=11; a=1}
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Formula = $row.a
    $ws.Range("A$r").NumberFormat = "0.00"
    $ws.Range("A$r").HorizontalAlignment = $xlRight

    $ws.Range("B$r").Formula = "=A$r*J`$3"
    $ws.Range("B$r").NumberFormat = "0.0"
    $ws.Range("B$r").HorizontalAlignment = $xlRight

    $ws.Range("C$r").Formula = "=J`$2*J`$2/B$r"
    $ws.Range("C$r").NumberFormat = "0.000"
    $ws.Range("C$r").HorizontalAlignment = $xlRight

    if ($r -eq 3) {
        $ws.Range("D$r").Value = 74.06
        $ws.Range("F$r").Value = 100.13
        $ws.Range("D$r").Interior.Color = 65535
        $ws.Range("F$r").Interior.Color = 65535
    } else {
        $ws.Range("D$r").Formula = "=D`$3"
        $ws.Range("F$r").Formula = "=F`$3"
    }
    $ws.Range("D$r").NumberFormat = "0.000"
    $ws.Range("D$r").HorizontalAlignment = $xlRight
    $ws.Range("F$r").NumberFormat = "0.000"
    $ws.Range("F$r").HorizontalAlignment = $xlRight

    $ws.Range("E$r").Formula = "=C$r-D$r"
    $ws.Range("E$r").NumberFormat = "0.000"

    $ws.Range("G$r").Formula = "=C$r-F$r"
    $ws.Range("G$r").NumberFormat = "0.000"

    $ws.Range("H$r").NumberFormat = "0.000"
}

# Bold label column (I) continues down, blank past row 3
foreach ($r in 4..11) {
    $ws.Range("I$r").Font.Bold = $true
}

# J6 - blank input cell referenced by the XS defined name, numFmt 0.00
$ws.Range("J6").NumberFormat = "0.00"

# ---- Row 12: SCR Xplant summary ----
$c = $ws.Range("D12"); $c.Value = "SCR Xplant"; $c.Font.Bold = $true
$c = $ws.Range("F12"); $c.Value = "SCR Xplant"; $c.Font.Bold = $true
$ws.Range("E12").Formula = "=`$J`$2*`$J`$2/`$D`$3/100"
$ws.Range("E12").NumberFormat = "0.000"
$ws.Range("G12").Formula = "=`$J`$2*`$J`$2/`$F`$3/100"
$ws.Range("G12").NumberFormat = "0.000"
$ws.Range("H12").NumberFormat = "0.000"

# ---- Column widths ----
$ws.Columns.Item(2).ColumnWidth = 9.36328125
$ws.Columns.Item(3).ColumnWidth = 9.36328125
$ws.Columns.Item(4).ColumnWidth = 9.36328125
$ws.Columns.Item(5).ColumnWidth = 9.54296875
$ws.Columns.Item(6).ColumnWidth = 9.54296875
$ws.Columns.Item(7).ColumnWidth = 9.54296875
$ws.Columns.Item(8).ColumnWidth = 9.54296875

# ---- Selection state ----
$ws.Range("A10").Select()

# ---- Weak sheet selection state ----
$wsWeak = $wb.Worksheets.Item("Weak")
$wsWeak.Range("F5").Select()
$ws.Select()

# ---- XS defined name now points at StepSCR!$J$6 ----
$n = $wb.Names.Item("StepSCR!XS")
$n.RefersTo = "=StepSCR!`$J`$6"
